# Update market/profit figures in the Cactuar_Profits leve-crafting workbook.
# Values below were refreshed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 155.96428
$ws.Range("I33").Value = 159.4
$ws.Range("K33").Value = 159.4
$ws.Range("M33").Value = 69.59999999999999
# Row 40
$ws.Range("H40").Value = 16697.334
$ws.Range("J40").Value = 15328.889
$ws.Range("L40").Value = 15328.889
$ws.Range("N40").Value = -15678.889
# Row 61
$ws.Range("H61").Value = 9999999
$ws.Range("I61").Value = 9999999
$ws.Range("K61").Value = 29999997
$ws.Range("M61").Value = -29999825
# Row 64
$ws.Range("H64").Value = 4095.1052
$ws.Range("I64").Value = 4023.5
$ws.Range("J64").Value = 4114.2
$ws.Range("K64").Value = 4023.5
$ws.Range("L64").Value = 4114.2
$ws.Range("M64").Value = -3775.5
$ws.Range("N64").Value = -4610.2
# Row 67
$ws.Range("H67").Value = 4095.1052
$ws.Range("I67").Value = 4023.5
$ws.Range("J67").Value = 4114.2
$ws.Range("K67").Value = 4023.5
$ws.Range("L67").Value = 4114.2
$ws.Range("M67").Value = -3165.5
$ws.Range("N67").Value = -5830.2
# Row 80
$ws.Range("H80").Value = 40001830
$ws.Range("I80").Value = 2254.4546
$ws.Range("J80").Value = 71430070
$ws.Range("K80").Value = 6763.3638
$ws.Range("L80").Value = 214290210
$ws.Range("M80").Value = -5765.3638
$ws.Range("N80").Value = -214292206
# Row 83
$ws.Range("H83").Value = 40001830
$ws.Range("I83").Value = 2254.4546
$ws.Range("J83").Value = 71430070
$ws.Range("K83").Value = 20290.0914
$ws.Range("L83").Value = 642870630
$ws.Range("M83").Value = -15298.0914
$ws.Range("N83").Value = -642880614
# Row 87
$ws.Range("H87").Value = 78247.5
$ws.Range("J87").Value = 78247.5
$ws.Range("L87").Value = 78247.5
$ws.Range("N87").Value = -80743.5
# Row 90
$ws.Range("H90").Value = 78247.5
$ws.Range("J90").Value = 78247.5
$ws.Range("L90").Value = 234742.5
$ws.Range("N90").Value = -247222.5
# Row 96
$ws.Range("H96").Value = 1879
$ws.Range("J96").Value = 3442.6667
$ws.Range("L96").Value = 10328.0001
$ws.Range("N96").Value = -13074.0001
# Row 100
$ws.Range("H100").Value = 1648.7059
$ws.Range("J100").Value = 1875
$ws.Range("L100").Value = 1875
$ws.Range("N100").Value = -2957
# Row 106
$ws.Range("H106").Value = 3613.3
$ws.Range("I106").Value = 3581.4443
$ws.Range("K106").Value = 3581.4443
$ws.Range("M106").Value = -2950.4443
# Row 112
$ws.Range("H112").Value = 3225.7297
$ws.Range("J112").Value = 3284.7778
$ws.Range("L112").Value = 9854.3334
$ws.Range("N112").Value = -12070.3334
# Row 132
$ws.Range("H132").Value = 3443.175
$ws.Range("I132").Value = 1635.2122
$ws.Range("J132").Value = 11966.429
$ws.Range("K132").Value = 4905.6366
$ws.Range("L132").Value = 35899.287
$ws.Range("M132").Value = -2375.6366
$ws.Range("N132").Value = -40959.287
# Row 137
$ws.Range("H137").Value = 12723685
$ws.Range("I137").Value = 672810.5600000001
$ws.Range("J137").Value = 27787278
$ws.Range("K137").Value = 2018431.68
$ws.Range("L137").Value = 83361834
$ws.Range("M137").Value = -2015881.68
$ws.Range("N137").Value = -83366934

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12403.966
$ws.Range("I32").Value = 11179.923
$ws.Range("J32").Value = 23012.334
$ws.Range("K32").Value = 11179.923
$ws.Range("L32").Value = 23012.334
$ws.Range("M32").Value = -10892.923
$ws.Range("N32").Value = -23586.334
# Row 45
$ws.Range("H45").Value = 2938.4211
$ws.Range("I45").Value = 2643.25
$ws.Range("K45").Value = 2643.25
$ws.Range("M45").Value = -2266.25
# Row 61
$ws.Range("H61").Value = 4117.3335
$ws.Range("I61").Value = 3963.6829
$ws.Range("J61").Value = 4448.8945
$ws.Range("K61").Value = 3963.6829
$ws.Range("L61").Value = 4448.8945
$ws.Range("M61").Value = -3751.6829
$ws.Range("N61").Value = -4872.8945
# Row 102
$ws.Range("H102").Value = 1720.2354
$ws.Range("I102").Value = 1749.5714
$ws.Range("K102").Value = 1749.5714
$ws.Range("M102").Value = -127.5714
# Row 132
$ws.Range("H132").Value = 3155.7595
$ws.Range("I132").Value = 1367.9672
$ws.Range("J132").Value = 9214.388999999999
$ws.Range("K132").Value = 4103.9016
$ws.Range("L132").Value = 27643.167
$ws.Range("M132").Value = -1573.9016
$ws.Range("N132").Value = -32703.167
# Row 136
$ws.Range("H136").Value = 4117.3335
$ws.Range("I136").Value = 3963.6829
$ws.Range("J136").Value = 4448.8945
$ws.Range("K136").Value = 11891.0487
$ws.Range("L136").Value = 13346.6835
$ws.Range("M136").Value = -9341.048699999999
$ws.Range("N136").Value = -18446.6835

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 999
$ws.Range("J14").Value = 999
$ws.Range("L14").Value = 999
$ws.Range("N14").Value = -1343
# Row 20
$ws.Range("H20").Value = 3781
$ws.Range("I20").Value = 3776.9092
$ws.Range("J20").Value = 3783.3684
$ws.Range("K20").Value = 3776.9092
$ws.Range("L20").Value = 3783.3684
$ws.Range("M20").Value = -3529.9092
$ws.Range("N20").Value = -4277.368399999999
# Row 99
$ws.Range("H99").Value = 4127.7
$ws.Range("I99").Value = 4696.7144
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 4696.7144
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -3198.7144
$ws.Range("N99").Value = -5796
# Row 105
$ws.Range("H105").Value = 2467.8572
$ws.Range("I105").Value = 2436.16
$ws.Range("K105").Value = 2436.16
$ws.Range("M105").Value = -689.1599999999999
# Row 107
$ws.Range("H107").Value = 5963.6665
$ws.Range("I107").Value = 4446.875
$ws.Range("K107").Value = 4446.875
$ws.Range("M107").Value = -2526.875
# Row 134
$ws.Range("H134").Value = 4868.6826
$ws.Range("I134").Value = 2320.7
$ws.Range("K134").Value = 6962.099999999999
$ws.Range("M134").Value = -4427.099999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 200
$ws.Range("I3").Value = 200
$ws.Range("K3").Value = 200
$ws.Range("M3").Value = -87
# Row 31
$ws.Range("H31").Value = 3063
$ws.Range("I31").Value = 1976.4517
$ws.Range("J31").Value = 9799.6
$ws.Range("K31").Value = 1976.4517
$ws.Range("L31").Value = 9799.6
$ws.Range("M31").Value = -1681.4517
$ws.Range("N31").Value = -10389.6
# Row 34
$ws.Range("H34").Value = 3063
$ws.Range("I34").Value = 1976.4517
$ws.Range("J34").Value = 9799.6
$ws.Range("K34").Value = 1976.4517
$ws.Range("L34").Value = 9799.6
$ws.Range("M34").Value = -1774.4517
$ws.Range("N34").Value = -10203.6
# Row 58
$ws.Range("H58").Value = 2663.9524
$ws.Range("I58").Value = 2639.7856
$ws.Range("J58").Value = 2712.2856
$ws.Range("K58").Value = 2639.7856
$ws.Range("L58").Value = 2712.2856
$ws.Range("M58").Value = -2436.7856
$ws.Range("N58").Value = -3118.2856
# Row 136
$ws.Range("H136").Value = 2663.9524
$ws.Range("I136").Value = 2639.7856
$ws.Range("J136").Value = 2712.2856
$ws.Range("K136").Value = 7919.3568
$ws.Range("L136").Value = 8136.8568
$ws.Range("M136").Value = -5369.3568
$ws.Range("N136").Value = -13236.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 4611490
$ws.Range("I4").Value = 252480.81
$ws.Range("K4").Value = 757442.4299999999
$ws.Range("M4").Value = -757330.4299999999
# Row 23
$ws.Range("H23").Value = 438.0909
$ws.Range("J23").Value = 451.9
$ws.Range("L23").Value = 1355.7
$ws.Range("N23").Value = -1825.7
# Row 131
$ws.Range("H131").Value = 18124714
$ws.Range("I131").Value = 9260997
$ws.Range("J131").Value = 22959470
$ws.Range("K131").Value = 27782991
$ws.Range("L131").Value = 68878410
$ws.Range("M131").Value = -27777951
$ws.Range("N131").Value = -68888490

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 18665.334
$ws.Range("J58").Value = 18665.334
$ws.Range("L58").Value = 18665.334
$ws.Range("N58").Value = -19219.334
# Row 97
$ws.Range("H97").Value = 1521.1364
$ws.Range("I97").Value = 1281.0588
$ws.Range("K97").Value = 1281.0588
$ws.Range("M97").Value = -785.0588
# Row 113
$ws.Range("H113").Value = 1663.8334
$ws.Range("I113").Value = 1359.3636
$ws.Range("J113").Value = 5013
$ws.Range("K113").Value = 1359.3636
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = 810.6364000000001
$ws.Range("N113").Value = -9353

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1681.0605
$ws.Range("I16").Value = 1260.6154
$ws.Range("J16").Value = 3242.7144
$ws.Range("K16").Value = 1260.6154
$ws.Range("L16").Value = 3242.7144
$ws.Range("M16").Value = -1090.6154
$ws.Range("N16").Value = -3582.7144
# Row 22
$ws.Range("H22").Value = 941.9474
$ws.Range("J22").Value = 1175.3334
$ws.Range("L22").Value = 1175.3334
$ws.Range("N22").Value = -1765.3334
# Row 27
$ws.Range("H27").Value = 941.9474
$ws.Range("J27").Value = 1175.3334
$ws.Range("L27").Value = 1175.3334
$ws.Range("N27").Value = -1389.3334
# Row 46
$ws.Range("H46").Value = 4403.6816
$ws.Range("I46").Value = 2998
$ws.Range("J46").Value = 4817.1177
$ws.Range("K46").Value = 2998
$ws.Range("L46").Value = 4817.1177
$ws.Range("M46").Value = -2810
$ws.Range("N46").Value = -5193.1177
# Row 93
$ws.Range("H93").Value = 2799.7
$ws.Range("I93").Value = 2116.4
$ws.Range("J93").Value = 3483
$ws.Range("K93").Value = 2116.4
$ws.Range("L93").Value = 3483
$ws.Range("M93").Value = -868.4000000000001
$ws.Range("N93").Value = -5979
# Row 102
$ws.Range("H102").Value = 119250
$ws.Range("J102").Value = 119250
$ws.Range("L102").Value = 119250
$ws.Range("N102").Value = -125740
# Row 132
$ws.Range("H132").Value = 3907.05
$ws.Range("I132").Value = 3258.0793
$ws.Range("J132").Value = 5012.054
$ws.Range("K132").Value = 9774.2379
$ws.Range("L132").Value = 15036.162
$ws.Range("M132").Value = -7244.2379
$ws.Range("N132").Value = -20096.162

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 808.6
$ws.Range("I107").Value = 613.1
$ws.Range("J107").Value = 1199.6
$ws.Range("K107").Value = 1839.3
$ws.Range("L107").Value = 3598.8
$ws.Range("M107").Value = 80.69999999999982
$ws.Range("N107").Value = -7438.799999999999
# Row 113
$ws.Range("H113").Value = 800.5625
$ws.Range("I113").Value = 391
$ws.Range("J113").Value = 2029.25
$ws.Range("K113").Value = 1173
$ws.Range("L113").Value = 6087.75
$ws.Range("M113").Value = 997
$ws.Range("N113").Value = -10427.75
# Row 132
$ws.Range("H132").Value = 6804373
$ws.Range("I132").Value = 12821274
$ws.Range("J132").Value = 2658
$ws.Range("K132").Value = 38463822
$ws.Range("L132").Value = 7974
$ws.Range("M132").Value = -38461292
$ws.Range("N132").Value = -13034
